$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fbln1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07741733333333332
$ws.Range("H2").Value = 0.232252
$ws.Range("I2").Value = 0.001631933472270046
$ws.Range("J2").Value = 0.001631933472270046
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 7.660191514149333
$ws.Range("R2").Value = 68.94172362734399
$ws.Range("S2").Value = 0.0003423812243137376
$ws.Range("T2").Value = 0.0003423812243137377

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fbln1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07741733333333332
$ws.Range("H3").Value = 0.232252
$ws.Range("I3").Value = 0.001631933472270046
$ws.Range("J3").Value = 0.001631933472270046
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 12.61950808201822
$ws.Range("R3").Value = 113.575572738164
$ws.Range("S3").Value = 0.0005640436821165194
$ws.Range("T3").Value = 0.0005640436821165196

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fbln1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07741733333333332
$ws.Range("H4").Value = 0.232252
$ws.Range("I4").Value = 0.001631933472270046
$ws.Range("J4").Value = 0.001631933472270046
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 5.062792188515556
$ws.Range("R4").Value = 45.56512969664001
$ws.Range("S4").Value = 0.0002262874217633028
$ws.Range("T4").Value = 0.0002262874217633029

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fbln1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.07741733333333332
$ws.Range("H5").Value = 0.232252
$ws.Range("I5").Value = 0.001631933472270046
$ws.Range("J5").Value = 0.001631933472270046
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 11.16921519047555
$ws.Range("R5").Value = 100.52293671428
$ws.Range("S5").Value = 0.0004992211440764857
$ws.Range("T5").Value = 0.0004992211440764858

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fbln1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 43.24729533333333
$ws.Range("H6").Value = 129.741886
$ws.Range("I6").Value = 0.9116396264352705
$ws.Range("J6").Value = 0.9116396264352705
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 4279.178195093821
$ws.Range("R6").Value = 38512.60375584439
$ws.Range("S6").Value = 0.1912628772774976
$ws.Range("T6").Value = 0.1912628772774977

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fbln1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 43.24729533333333
$ws.Range("H7").Value = 129.741886
$ws.Range("I7").Value = 0.9116396264352705
$ws.Range("J7").Value = 0.9116396264352705
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 7049.578815051267
$ws.Range("R7").Value = 63446.2093354614
$ws.Range("S7").Value = 0.3150891751381332
$ws.Range("T7").Value = 0.3150891751381332

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fbln1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 43.24729533333333
$ws.Range("H8").Value = 129.741886
$ws.Range("I8").Value = 0.9116396264352705
$ws.Range("J8").Value = 0.9116396264352705
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 2828.204738663503
$ws.Range("R8").Value = 25453.84264797152
$ws.Range("S8").Value = 0.126409920593357
$ws.Range("T8").Value = 0.1264099205933571

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fbln1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 43.24729533333333
$ws.Range("H9").Value = 129.741886
$ws.Range("I9").Value = 0.9116396264352705
$ws.Range("J9").Value = 0.9116396264352705
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 6239.408246009282
$ws.Range("R9").Value = 56154.67421408354
$ws.Range("S9").Value = 0.2788776534262826
$ws.Range("T9").Value = 0.2788776534262826

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fbln1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.114312666666667
$ws.Range("H10").Value = 12.342938
$ws.Range("I10").Value = 0.08672844009245946
$ws.Range("J10").Value = 0.08672844009245947
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 407.0977598783707
$ws.Range("R10").Value = 3663.879838905336
$ws.Range("S10").Value = 0.01819571079718821
$ws.Range("T10").Value = 0.01819571079718822

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Fbln1"
$ws.Range("C11").Value = "Itgb1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.114312666666667
$ws.Range("H11").Value = 12.342938
$ws.Range("I11").Value = 0.08672844009245946
$ws.Range("J11").Value = 0.08672844009245947
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 670.6586201490185
$ws.Range("R11").Value = 6035.927581341166
$ws.Range("S11").Value = 0.02997587188767334
$ws.Range("T11").Value = 0.02997587188767335

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Fbln1"
$ws.Range("C12").Value = "Itgb1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.114312666666667
$ws.Range("H12").Value = 12.342938
$ws.Range("I12").Value = 0.08672844009245946
$ws.Range("J12").Value = 0.08672844009245947
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 269.0600300093512
$ws.Range("R12").Value = 2421.540270084161
$ws.Range("S12").Value = 0.01202595291753913
$ws.Range("T12").Value = 0.01202595291753913

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Fbln1"
$ws.Range("C13").Value = "Itgb1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.114312666666667
$ws.Range("H13").Value = 12.342938
$ws.Range("I13").Value = 0.08672844009245946
$ws.Range("J13").Value = 0.08672844009245947
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 593.5833947810911
$ws.Range("R13").Value = 5342.25055302982
$ws.Range("S13").Value = 0.02653090449005877
$ws.Range("T13").Value = 0.02653090449005878
